$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 115.25
$ws.Range("I33").Value = 113.8
$ws.Range("J33").Value = 116.28571
$ws.Range("K33").Value = 113.8
$ws.Range("L33").Value = 116.28571
$ws.Range("M33").Value = 115.2
$ws.Range("N33").Value = -574.28571

$ws.Range("H88").Value = 2128.7
$ws.Range("I88").Value = 5496
$ws.Range("J88").Value = 1286.875
$ws.Range("K88").Value = 5496
$ws.Range("L88").Value = 1286.875
$ws.Range("M88").Value = -5090
$ws.Range("N88").Value = -2098.875

$ws.Range("H91").Value = 2128.7
$ws.Range("I91").Value = 5496
$ws.Range("J91").Value = 1286.875
$ws.Range("K91").Value = 5496
$ws.Range("L91").Value = 1286.875
$ws.Range("M91").Value = -4092
$ws.Range("N91").Value = -4094.875

$ws.Range("H98").Value = 1933.5652
$ws.Range("I98").Value = 1976
$ws.Range("K98").Value = 1976
$ws.Range("M98").Value = -478

$ws.Range("H100").Value = 1474.3636
$ws.Range("I100").Value = 1321.8
$ws.Range("K100").Value = 1321.8
$ws.Range("M100").Value = -780.8

$ws.Range("H107").Value = 1515
$ws.Range("I107").Value = 966
$ws.Range("J107").Value = 2338.5
$ws.Range("K107").Value = 966
$ws.Range("L107").Value = 2338.5
$ws.Range("M107").Value = 954
$ws.Range("N107").Value = -6178.5

$ws.Range("H116").Value = 16159.9
$ws.Range("J116").Value = 7685.5713
$ws.Range("L116").Value = 7685.5713
$ws.Range("N116").Value = -14569.5713

$ws.Range("H122").Value = 1933.5652
$ws.Range("I122").Value = 1976
$ws.Range("K122").Value = 5928
$ws.Range("M122").Value = -3478

$ws.Range("H125").Value = 459.7143
$ws.Range("I125").Value = 369.66666
$ws.Range("K125").Value = 3326.99994
$ws.Range("M125").Value = -866.9999399999997

$ws.Range("H132").Value = 1118.5758
$ws.Range("I132").Value = 1072.862
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 3218.586
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -688.5860000000002
$ws.Range("N132").Value = -9410

$ws.Range("H137").Value = 1420.7778
$ws.Range("I137").Value = 1337.4
$ws.Range("K137").Value = 4012.2
$ws.Range("M137").Value = -1462.2

$ws.Range("H138").Value = 2889.68
$ws.Range("J138").Value = 2841.6
$ws.Range("L138").Value = 8524.799999999999
$ws.Range("N138").Value = -18804.8

$ws.Range("H141").Value = 1477133.4
$ws.Range("I141").Value = 2156265.2
$ws.Range("K141").Value = 6468795.600000001
$ws.Range("M141").Value = -6463615.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4270.8696
$ws.Range("I32").Value = 3007.2632
$ws.Range("K32").Value = 3007.2632
$ws.Range("M32").Value = -2720.2632

$ws.Range("H62").Value = 29999
$ws.Range("J62").Value = 29999
$ws.Range("L62").Value = 29999
$ws.Range("N62").Value = -31247

$ws.Range("H65").Value = 29999
$ws.Range("J65").Value = 29999
$ws.Range("L65").Value = 89997
$ws.Range("N65").Value = -96237

$ws.Range("H74").Value = 1105.3673
$ws.Range("I74").Value = 746.1951
$ws.Range("J74").Value = 2946.125
$ws.Range("K74").Value = 746.1951
$ws.Range("L74").Value = 2946.125
$ws.Range("M74").Value = 127.8049
$ws.Range("N74").Value = -4694.125

$ws.Range("H77").Value = 1105.3673
$ws.Range("I77").Value = 746.1951
$ws.Range("J77").Value = 2946.125
$ws.Range("K77").Value = 3730.9755
$ws.Range("L77").Value = 14730.625
$ws.Range("M77").Value = 637.0245
$ws.Range("N77").Value = -23466.625

$ws.Range("H110").Value = 2654.8
$ws.Range("I110").Value = 1337
$ws.Range("K110").Value = 1337
$ws.Range("M110").Value = 708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I86").Value = 2312.6365
$ws.Range("J86").Value = 335916.5
$ws.Range("K86").Value = 2312.6365
$ws.Range("L86").Value = 335916.5
$ws.Range("M86").Value = -1189.6365
$ws.Range("N86").Value = -338162.5

$ws.Range("I89").Value = 2312.6365
$ws.Range("J89").Value = 335916.5
$ws.Range("K89").Value = 11563.1825
$ws.Range("L89").Value = 1679582.5
$ws.Range("M89").Value = -5947.182500000001
$ws.Range("N89").Value = -1690814.5

$ws.Range("H94").Value = 907.8570999999999
$ws.Range("I94").Value = 907.8570999999999
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 907.8570999999999
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -456.8570999999999
$ws.Range("N94").ClearContents()

$ws.Range("H99").Value = 1504.4286
$ws.Range("I99").Value = 1278.3636
$ws.Range("K99").Value = 1278.3636
$ws.Range("M99").Value = 219.6364000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 319.9
$ws.Range("I7").Value = 183.16667
$ws.Range("K7").Value = 183.16667
$ws.Range("M7").Value = -70.16667000000001

$ws.Range("H31").Value = 2670.762
$ws.Range("J31").Value = 3799.5293
$ws.Range("L31").Value = 3799.5293
$ws.Range("N31").Value = -4389.5293

$ws.Range("H34").Value = 2670.762
$ws.Range("J34").Value = 3799.5293
$ws.Range("L34").Value = 3799.5293
$ws.Range("N34").Value = -4203.5293

$ws.Range("H99").Value = 2041.3
$ws.Range("I99").Value = 1249.8334
$ws.Range("J99").Value = 3228.5
$ws.Range("K99").Value = 1249.8334
$ws.Range("L99").Value = 3228.5
$ws.Range("M99").Value = 248.1666
$ws.Range("N99").Value = -6224.5

$ws.Range("H126").Value = 2041.3
$ws.Range("I126").Value = 1249.8334
$ws.Range("J126").Value = 3228.5
$ws.Range("K126").Value = 3749.5002
$ws.Range("L126").Value = 9685.5
$ws.Range("M126").Value = -1279.5002
$ws.Range("N126").Value = -14625.5

$ws.Range("H132").Value = 2097.2693
$ws.Range("I132").Value = 1151
$ws.Range("K132").Value = 3453
$ws.Range("M132").Value = -923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 924.75
$ws.Range("I11").Value = 605.5
$ws.Range("J11").Value = 1244
$ws.Range("K11").Value = 1816.5
$ws.Range("L11").Value = 3732
$ws.Range("M11").Value = -1676.5
$ws.Range("N11").Value = -4012

$ws.Range("H98").Value = 789.8
$ws.Range("J98").Value = 937.25
$ws.Range("L98").Value = 2811.75
$ws.Range("N98").Value = -5807.75

$ws.Range("H131").Value = 10067.465
$ws.Range("J131").Value = 11176.961
$ws.Range("L131").Value = 33530.883
$ws.Range("N131").Value = -43610.883

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 43666
$ws.Range("I19").Value = 21659.334
$ws.Range("K19").Value = 21659.334
$ws.Range("M19").Value = -21371.334

$ws.Range("H127").Value = 13995
$ws.Range("J127").Value = 13995
$ws.Range("L127").Value = 13995
$ws.Range("N127").Value = -23915

$ws.Range("H139").Value = 43199.6
$ws.Range("J139").Value = 43199.6
$ws.Range("L139").Value = 43199.6
$ws.Range("N139").Value = -53479.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5333.8335
$ws.Range("I7").Value = 4001
$ws.Range("K7").Value = 4001
$ws.Range("M7").Value = -3889

$ws.Range("H16").Value = 3215.1177
$ws.Range("I16").Value = 3740.111
$ws.Range("K16").Value = 3740.111
$ws.Range("M16").Value = -3570.111

$ws.Range("H126").Value = 5333.8335
$ws.Range("I126").Value = 4001
$ws.Range("K126").Value = 12003
$ws.Range("M126").Value = -9533

$ws.Range("H136").Value = 3488.25
$ws.Range("I136").Value = 1640.8889
$ws.Range("J136").Value = 4999.727
$ws.Range("K136").Value = 4922.6667
$ws.Range("L136").Value = 14999.181
$ws.Range("M136").Value = -2372.6667
$ws.Range("N136").Value = -20099.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1429.5454
$ws.Range("I81").Value = 1155.625
$ws.Range("K81").Value = 2311.25
$ws.Range("M81").Value = -1250.25

$ws.Range("H84").Value = 1429.5454
$ws.Range("I84").Value = 1155.625
$ws.Range("K84").Value = 11556.25
$ws.Range("M84").Value = -6252.25

$ws.Range("H132").Value = 3277.2144
$ws.Range("I132").Value = 2654.4443
$ws.Range("J132").Value = 4398.2
$ws.Range("K132").Value = 7963.3329
$ws.Range("L132").Value = 13194.6
$ws.Range("M132").Value = -5433.3329
$ws.Range("N132").Value = -18254.6
